# Add a new worksheet "Hoja3" after "Hoja2", give it the value "hide" in
# A1, and make it the active (selected) sheet.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the last existing sheet (Hoja2), so it lands in
# the 3rd tab position.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Hoja3"

# Populate its A1 cell.
$newSheet.Range("A1").Value = "hide"

# Make the new sheet the active tab.
$newSheet.Activate() | Out-Null
$newSheet.Range("A2").Select() | Out-Null
